$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Primärdaten")

# Remember column N's width so the freshly inserted column O can match it
# (mirrors Excel's native "insert column" behaviour of inheriting the
# left-neighbour's width).
$nWidth = $ws.Columns.Item(14).ColumnWidth

# Insert a new column before column O (15) - shifts old O (and beyond) right.
$ws.Columns.Item(15).Insert()
$ws.Columns.Item(15).ColumnWidth = $nWidth

# New header for the inserted column.
$ws.Cells.Item(1, 15).Value = "DEPENDENCY"

# Expand the _FilterDatabase defined name to include the new column.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Primärdaten!`$D`$2:`$P`$58"
    }
}

# Update view: scroll so column H is the left-most visible column, and select O1.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("O1").Select() | Out-Null
